# Updated symbol list on Sun Jan 29 20:22:07 UTC 2023 with GitHub Actions
#
# This refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# ranking sheet with newly scraped quotes, and reflects that rank #46/#47
# (CoinbaseStockToken and BOLO) swapped positions in the source ranking on
# this run.
#
# All values in this sheet are stored as literal text (e.g. "317.00",
# "3.30%") rather than numbers, so every assignment below is prefixed with
# a leading apostrophe. That is Excel's standard "treat as text" marker:
# it forces the cell to stay a text value with the exact characters typed,
# instead of being auto-parsed into a number/percentage - matching the
# original workbook's text-cell formatting and avoiding unwanted numeric
# reformatting (e.g. "317.81" staying "317.81" rather than becoming 317.81
# as a General-formatted number, and "3.64%" staying literal text instead
# of the numeric fraction 0.0364).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'317.81"
$ws.Range('E2').Value = "'3.64%"
$ws.Range('D3').Value = "'39.89"
$ws.Range('E3').Value = "'2.58%"
$ws.Range('D4').Value = "'5.143"
$ws.Range('E4').Value = "'0.74%"
$ws.Range('D5').Value = "'0.08225"
$ws.Range('E5').Value = "'1.70%"
$ws.Range('D6').Value = "'2.067"
$ws.Range('E6').Value = "'6.43%"
$ws.Range('D7').Value = "'8.351"
$ws.Range('E7').Value = "'4.65%"
$ws.Range('D8').Value = "'4.335"
$ws.Range('E8').Value = "'3.64%"
$ws.Range('D9').Value = "'0.9373"
$ws.Range('E9').Value = "'0.69%"
$ws.Range('E10').Value = "'-7.13%"
$ws.Range('D11').Value = "'0.1984"
$ws.Range('E11').Value = "'2.83%"
$ws.Range('D12').Value = "'0.09129"
$ws.Range('E12').Value = "'-0.48%"
$ws.Range('D13').Value = "'0.03512"
$ws.Range('E13').Value = "'0.26%"
$ws.Range('D14').Value = "'0.09812"
$ws.Range('E14').Value = "'0.34%"
$ws.Range('D15').Value = "'0.001407"
$ws.Range('E15').Value = "'1.21%"
$ws.Range('D16').Value = "'0.006353"
$ws.Range('E16').Value = "'6.31%"
$ws.Range('D17').Value = "'3.698"
$ws.Range('E17').Value = "'-2.18%"
$ws.Range('D18').Value = "'3.206"
$ws.Range('E18').Value = "'-6.42%"
$ws.Range('D19').Value = "'0.3499"
$ws.Range('E19').Value = "'1.14%"
$ws.Range('D20').Value = "'0.1315"
$ws.Range('E20').Value = "'1.02%"
$ws.Range('D21').Value = "'4.962"
$ws.Range('E21').Value = "'5.94%"
$ws.Range('D22').Value = "'0.2451"
$ws.Range('E22').Value = "'1.55%"
$ws.Range('D23').Value = "'0.04355"
$ws.Range('E23').Value = "'-0.40%"
$ws.Range('D24').Value = "'0.001229"
$ws.Range('E24').Value = "'-0.61%"
$ws.Range('D25').Value = "'0.004826"
$ws.Range('E25').Value = "'12.71%"
$ws.Range('D26').Value = "'0.0001298"
$ws.Range('E26').Value = "'-0.28%"
$ws.Range('D27').Value = "'0.0004000"
$ws.Range('E27').Value = "'-10.06%"
$ws.Range('D39').Value = "'0.02209"
$ws.Range('E39').Value = "'8.35%"
$ws.Range('D40').Value = "'0.05223"
$ws.Range('E40').Value = "'3.06%"
$ws.Range('D41').Value = "'0.007736"
$ws.Range('E41').Value = "'2.58%"
$ws.Range('D42').Value = "'0.009696"
$ws.Range('E42').Value = "'-5.44%"
$ws.Range('D43').Value = "'0.1410"
$ws.Range('E43').Value = "'4.50%"
$ws.Range('D44').Value = "'0.002047"
$ws.Range('E44').Value = "'-3.62%"
$ws.Range('D45').Value = "'0.009648"
$ws.Range('E45').Value = "'-2.75%"
$ws.Range('D46').Value = "'0.00006622"
$ws.Range('E46').Value = "'6.89%"
$ws.Range('D47').Value = "'0.00000000750"
$ws.Range('E47').Value = "'-0.15%"
$ws.Range('B48').Value = "'BOLO"
$ws.Range('C48').Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range('D48').Value = "'0.002939"
$ws.Range('E48').Value = "'-5.52%"
$ws.Range('B49').Value = "'CoinbaseStockToken"
$ws.Range('C49').Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range('D49').Value = "'0.001690"
$ws.Range('E49').Value = "'5.53%"
$ws.Range('D50').Value = "'0.00002100"
$ws.Range('E50').Value = "'-0.15%"
$ws.Range('D51').Value = "'0.0002000"
$ws.Range('E51').Value = "'-0.15%"
